$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201 (shifts existing rows 201-271 down to 202-272)
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new data record
$ws.Cells.Item(201, 1).Value = 4
$ws.Cells.Item(201, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(201, 3).Value = "Los Lagos"
$ws.Cells.Item(201, 4).Value = 44559
$ws.Cells.Item(201, 5).Value = 10
$ws.Cells.Item(201, 6).Value = 100114013
$ws.Cells.Item(201, 7).Value = "Zanahoria"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 300
$ws.Cells.Item(201, 11).Value = 11500
$ws.Cells.Item(201, 12).Value = 12000
$ws.Cells.Item(201, 13).Value = 11750
$ws.Cells.Item(201, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(201, 15).Value = "Región de Ñuble"
$ws.Cells.Item(201, 16).Value = 588
$ws.Cells.Item(201, 17).Value = 20
$ws.Cells.Item(201, 18).Value = "Hortaliza"
